$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 entirely (shifts everything below up by 2 rows,
# and adjusts shared formula ranges automatically).
$ws.Range("A3:A4").EntireRow.Delete()

# Reset the frozen-pane scroll position back to the top of the data.
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 5

# Select entire rows 3:4 (the rows that slid up into the deleted rows'
# place), matching the post-delete selection state left in the file.
$ws.Range("A3:A4").EntireRow.Select()
